# Update the cryptos price/volume table with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume }
# (only columns that actually changed are present per row)
$updates = @{
    2  = @{ D = "27.959.57"; E = "  -0.08%  " }
    3  = @{ D = "1.638.39";  E = "  -0.58%  " }
    4  = @{ E = "  +0.07%  " }
    5  = @{ D = "212.36";    E = "  -0.77%  " }
    6  = @{ E = "  -0.07%  " }
    7  = @{ E = "  +0.05%  " }
    8  = @{ D = "23.33";     E = "  -1.57%  " }
    9  = @{ E = "  -2.53%  " }
    10 = @{ E = "  -0.11%  " }
    11 = @{ D = "0.0881";    E = "  +0.94%  " }
    12 = @{ D = "1.870.91";  E = "  -0.50%  " }
    13 = @{ D = "1.653.20";  E = "  +0.35%  " }
    14 = @{ E = "  -0.24%  " }
    15 = @{ E = "  +0.76%  " }
    16 = @{ D = "65.37";     E = "  -0.64%  " }
    17 = @{ D = "27.964.34"; E = "  +0.00%  " }
    18 = @{ D = "231.37";    E = "  -0.54%  " }
    19 = @{ D = "0.0₃0722";  E = "  -0.22%  " }
    20 = @{ E = "  -1.85%  " }
    21 = @{ E = "  +0.08%  " }
    22 = @{ E = "  -2.36%  " }
    23 = @{ D = "4.37";      E = "  -0.51%  " }
    24 = @{ E = "  -4.14%  " }
    25 = @{ D = "153.57";    E = "  +1.37%  " }
    26 = @{ D = "6.99";      E = "  +0.68%  " }
    27 = @{ E = "  -0.38%  " }
    28 = @{ D = "15.62";     E = "  -0.76%  " }
    30 = @{ D = "1.19";      E = "  -0.71%  " }
    31 = @{ E = "  -0.25%  " }
    32 = @{ E = "  +1.63%  " }
    33 = @{ D = "1.407.36";  E = "  -3.36%  " }
    34 = @{ E = "  -1.50%  " }
    35 = @{ E = "  +1.14%  " }
    36 = @{ E = "  +1.59%  " }
    37 = @{ E = "  +0.27%  " }
    38 = @{ E = "  -0.35%  " }
    39 = @{ D = "0.927";     E = "  +1.10%  " }
    40 = @{ E = "  -1.46%  " }
    41 = @{ E = "  +0.20%  " }
    42 = @{ E = "  -0.02%  " }
    43 = @{ D = "67.05" }
    44 = @{ E = "  +2.54%  " }
    45 = @{ E = "  +1.45%  " }
    46 = @{ E = "  -1.02%  " }
    47 = @{ D = "1.780.14";  E = "  -0.56%  " }
    48 = @{ D = "87.99";     E = "  -0.94%  " }
    49 = @{ E = "  -0.76%  " }
    50 = @{ E = "  -0.50%  " }
    51 = @{ D = "7.57";      E = "  -2.45%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        # Force text storage so numeric-looking values (e.g. "212.36") are
        # not auto-converted into floating point numbers by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData.D
        $cell.ClearFormats()
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$row").Value = $rowData.E
    }
}
